$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Data row template shared by rows 45 and 46 (columns B:O)
$rowValues = @(8, 6, 211, 386, 362, 388, 2681, 388, 1216, 119, 304, 30, 3087, 4051)

# Row 45
$ws.Cells.Item(45, 1).Value = 45700.60665509259
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(45, $i + 2).Value = $rowValues[$i]
}

# Row 46
$ws.Cells.Item(46, 1).Value = 45700.615081018521
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(46, $i + 2).Value = $rowValues[$i]
}

# Copy formatting (number format / style) from the previous data row (44) to
# the two new rows, so the new cells keep the same styles (s="5"/s="3").
$srcRow = $ws.Range("A44:O44")
$dstRow1 = $ws.Range("A45:O45")
$dstRow2 = $ws.Range("A46:O46")
$srcRow.Copy()
$dstRow1.PasteSpecial(-4122)  # xlPasteFormats
$dstRow2.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
